$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated how total kelp cover is calculated: the "Mussel Point" row
# (row 25) now counts as entered rather than not entered.
$ws.Cells.Item(25, 6).Value = "entered"

# New campaign rows for Pangnirtung NFA survey sites.
$ws.Cells.Item(31, 1).Value = "Pangnirtung NFA"
$ws.Cells.Item(31, 2).Value = 43746
$ws.Cells.Item(31, 3).Value = "Pangnirtung 1"
$ws.Cells.Item(31, 4).Formula = "=66+(6.4/60)"
$ws.Cells.Item(31, 5).Formula = "=-59.5/60-65"
$ws.Cells.Item(31, 6).Value = "entered"

$ws.Cells.Item(32, 1).Value = "Pangnirtung NFA"
$ws.Cells.Item(32, 2).Value = 43746
$ws.Cells.Item(32, 3).Value = "Pangnirtung 2"
$ws.Cells.Item(32, 4).Formula = "=66+(15.96/60)"
$ws.Cells.Item(32, 5).Formula = "=-59.5/60-67"
$ws.Cells.Item(32, 6).Value = "not entered"

# Match the date formatting used by the other rows in column B.
$ws.Cells.Item(30, 2).Copy() | Out-Null
$ws.Cells.Item(31, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(32, 2).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the selection to reflect where the author ended up after entering data.
$ws.Range("F33").Select() | Out-Null
